$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "63.165.20"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "3.152.27"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'589.89"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").Value = "'138.04"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.148.03"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'5.30"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").Value = "'34.16"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "3.673.71"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "3.152.11"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "63.130.15"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "'6.67"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").Value = "'476.35"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'14.04"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "'84.67"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("D25").Value = "'12.99"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "'7.97"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'26.97"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").Value = "'52.77"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'5.81"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").Value = "0.0₃0698"
$ws.Range("E38").Value = "  -7.06%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "'419.34"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  -7.68%  "
$ws.Range("D42").Value = "'8.28"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").Value = "2.931.13"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "'25.46"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -8.59%  "
$ws.Range("E51").Value = "  -0.24%  "
